$wb = $excel.ActiveWorkbook

# --- Metadata sheet updates ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value = "2025-06-13T15:45:04+00:00"
$meta.Range("B15").Value = "4.0.1"

# --- Elements sheet updates ---
$el = $wb.Worksheets.Item("Elements")

# Extension.id Type(s): "id\n" -> "string\n"
$el.Range("K3").Value = "string`n"

# Root Extension Constraint(s): simplify the ele-1 rule (drop the Parameters clause)
$el.Range("AJ2").Value = "ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}`next-1:Must have either extensions or value[x], not both {extension.exists() != value.exists()}"

# Extension.value[x] Definition: R4B -> R4 in the extensibility link
$el.Range("M6").Value = "Value of extension - must be one of a constrained set of the data types (see [Extensibility](http://hl7.org/fhir/R4/extensibility.html) for a list)."
